# Remove the "population_density_per_sq_km" column (column M) from Sheet1,
# including its header and all data rows, shifting the remaining columns
# left. The companion "raw" sheet is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns("M:M").Delete()

# Reflect the cursor's resting position after the edit, matching where the
# author last clicked in the workbook.
[void]$ws.Range("S16").Select()
